# XF-980 AUTO_TC 6.4.8 Roles - Create new Role - Create new role
# Adds a "Role Name" / "Role Permissions" column pair to the "6_Tenants" sheet
# and fills in sample data for the newly created QA test role.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("6_Tenants")

# New header cells (row 1) following the existing "BK Color" header in AZ1.
$ws.Range("BB1").Value = "Role Permissions"
$ws.Range("BA1").Value = "Role Name"

# New data cells (row 2) describing the newly created QA test role.
$ws.Range("BB2").Value = "Tenant Delete"
$ws.Range("BA2").Value = "QATest Role "

# Match the author's final selection/view state on this sheet.
$ws.Activate()
$ws.Range("BC9").Select()
